# final_ass_time_changed/data/losses_tables/A.1_lossestable.xlsx
# "aanpassing tabellen voor extra timesteps"
#
# Inserts two new interpolated rows (timesteps) into the losses table,
# located between the existing timesteps, so that the table goes from
# timesteps 0,1,2 to timesteps 0,1,2,3,4 where the new odd timesteps
# (1 and 3) are the average of their neighbours, and the old timesteps
# 1 and 2 are shifted down to become timesteps 2 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above the old row 4 (old A=2); it becomes new row 5 ---
$ws.Rows.Item(4).Insert()

# --- Insert a new row above the old row 3 (old A=1); it becomes new row 3 ---
$ws.Rows.Item(3).Insert()

# After the two inserts, the rows now look like:
#   row1 = header
#   row2 = old row2 (A=0) unchanged
#   row3 = blank (new)
#   row4 = old row3 (A=1), now holds the values for A=2
#   row5 = blank (new)
#   row6 = old row4 (A=2), now holds the values for A=4

# --- Copy the "timestep" column (A) formatting onto the new rows ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# --- New row 3 = average of row 2 and row 4 ---
$ws.Range("B3").Formula = "=(B2+B4)/2"
$ws.Range("C3").Formula = "=(C2+C4)/2"
$ws.Range("D3").Formula = "=(D2+D4)/2"
$ws.Range("E3").Formula = "=(E2+E4)/2"
$ws.Range("F3").Formula = "=(F2+F4)/2"
$ws.Range("G3").Formula = "=(G2+G4)/2"
$ws.Range("H3").Formula = "=(H2+H4)/2"

# --- New row 5 = average of row 4 and row 6 ---
$ws.Range("B5").Formula = "=(B4+B6)/2"
$ws.Range("C5").Formula = "=(C4+C6)/2"
$ws.Range("D5").Formula = "=(D4+D6)/2"
$ws.Range("E5").Formula = "=(E4+E6)/2"
$ws.Range("F5").Formula = "=(F4+F6)/2"
$ws.Range("G5").Formula = "=(G4+G6)/2"
$ws.Range("H5").Formula = "=(H4+H6)/2"

# New rows B..H should use the plain (unstyled) formatting like row2/row4/row6
# - clear the inherited border/bold style the Insert() copied down from the
# row above.
$ws.Range("B3:H3").Style = "Normal"
$ws.Range("B5:H5").Style = "Normal"

# --- Update the selection to match the saved workbook's cursor position ---
$ws.Range("G9").Select() | Out-Null

$wb.Save()
